$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1725.75
$ws.Range("J62").Value = 1866
$ws.Range("L62").Value = 1866
$ws.Range("N62").Value = -3114
$ws.Range("H65").Value = 1725.75
$ws.Range("J65").Value = 1866
$ws.Range("L65").Value = 9330
$ws.Range("N65").Value = -15570
$ws.Range("H106").Value = 16011.429
$ws.Range("I106").Value = 18283.334
$ws.Range("K106").Value = 18283.334
$ws.Range("M106").Value = -17652.334
$ws.Range("H129").Value = 2633.3333
$ws.Range("I129").Value = 666.6667
$ws.Range("K129").Value = 2000.0001
$ws.Range("M129").Value = 2999.9999
$ws.Range("H132").Value = 3970542
$ws.Range("I132").Value = 2038.6897
$ws.Range("J132").Value = 50005180
$ws.Range("K132").Value = 6116.0691
$ws.Range("L132").Value = 150015540
$ws.Range("M132").Value = -3586.0691
$ws.Range("N132").Value = -150020600
$ws.Range("H135").Value = 1036.1489
$ws.Range("I135").Value = 969.3570999999999
$ws.Range("J135").Value = 1597.2
$ws.Range("K135").Value = 8724.213899999999
$ws.Range("L135").Value = 14374.8
$ws.Range("M135").Value = -6189.213899999999
$ws.Range("N135").Value = -19444.8
$ws.Range("H137").Value = 48001.332
$ws.Range("I137").Value = 3501
$ws.Range("J137").Value = 70251.5
$ws.Range("K137").Value = 10503
$ws.Range("L137").Value = 210754.5
$ws.Range("M137").Value = -7953
$ws.Range("N137").Value = -215854.5
$ws.Range("H138").Value = 3761.57
$ws.Range("I138").Value = 1700.5555
$ws.Range("J138").Value = 4920.8906
$ws.Range("K138").Value = 5101.666499999999
$ws.Range("L138").Value = 14762.6718
$ws.Range("M138").Value = 38.33350000000064
$ws.Range("N138").Value = -25042.6718

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18223.408
$ws.Range("I32").Value = 18189.604
$ws.Range("J32").Value = 18933.334
$ws.Range("K32").Value = 18189.604
$ws.Range("L32").Value = 18933.334
$ws.Range("M32").Value = -17902.604
$ws.Range("N32").Value = -19507.334
$ws.Range("H119").Value = 60000
$ws.Range("J119").Value = 60000
$ws.Range("L119").Value = 60000
$ws.Range("N119").Value = -69676
$ws.Range("H132").Value = 2750.5417
$ws.Range("I132").Value = 2349.4119
$ws.Range("J132").Value = 3724.7144
$ws.Range("K132").Value = 7048.2357
$ws.Range("L132").Value = 11174.1432
$ws.Range("M132").Value = -4518.2357
$ws.Range("N132").Value = -16234.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 89965
$ws.Range("J59").Value = 89965
$ws.Range("L59").Value = 89965
$ws.Range("N59").Value = -91659
$ws.Range("H140").Value = 48780
$ws.Range("J140").Value = 48780
$ws.Range("L140").Value = 48780
$ws.Range("N140").Value = -59140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 43790
$ws.Range("J52").Value = 43790
$ws.Range("L52").Value = 43790
$ws.Range("N52").Value = -44378
$ws.Range("H105").Value = 2475.3333
$ws.Range("I105").Value = 2522.889
$ws.Range("J105").Value = 2332.6667
$ws.Range("K105").Value = 2522.889
$ws.Range("L105").Value = 2332.6667
$ws.Range("M105").Value = -775.8890000000001
$ws.Range("N105").Value = -5826.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1017
$ws.Range("I98").Value = 240
$ws.Range("J98").Value = 1572
$ws.Range("K98").Value = 720
$ws.Range("L98").Value = 4716
$ws.Range("M98").Value = 778
$ws.Range("N98").Value = -7712
$ws.Range("H107").Value = 455.81818
$ws.Range("I107").Value = 190.75
$ws.Range("J107").Value = 607.2857
$ws.Range("K107").Value = 572.25
$ws.Range("L107").Value = 1821.8571
$ws.Range("M107").Value = 1347.75
$ws.Range("N107").Value = -5661.8571
$ws.Range("H113").Value = 2349.6924
$ws.Range("I113").Value = 2267.6667
$ws.Range("J113").Value = 2420
$ws.Range("K113").Value = 6803.000100000001
$ws.Range("L113").Value = 7260
$ws.Range("M113").Value = -4633.000100000001
$ws.Range("N113").Value = -11600

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2695.2307
$ws.Range("I80").Value = 2520.2
$ws.Range("J80").Value = 2933.9092
$ws.Range("K80").Value = 2520.2
$ws.Range("L80").Value = 2933.9092
$ws.Range("M80").Value = -1522.2
$ws.Range("N80").Value = -4929.9092
$ws.Range("H83").Value = 2695.2307
$ws.Range("I83").Value = 2520.2
$ws.Range("J83").Value = 2933.9092
$ws.Range("K83").Value = 12601
$ws.Range("L83").Value = 14669.546
$ws.Range("M83").Value = -7609
$ws.Range("N83").Value = -24653.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1423.619
$ws.Range("I22").Value = 750
$ws.Range("J22").Value = 1582.1177
$ws.Range("K22").Value = 750
$ws.Range("L22").Value = 1582.1177
$ws.Range("M22").Value = -455
$ws.Range("N22").Value = -2172.1177
$ws.Range("H27").Value = 1423.619
$ws.Range("I27").Value = 750
$ws.Range("J27").Value = 1582.1177
$ws.Range("K27").Value = 750
$ws.Range("L27").Value = 1582.1177
$ws.Range("M27").Value = -643
$ws.Range("N27").Value = -1796.1177
$ws.Range("H61").Value = 1743.2307
$ws.Range("I61").Value = 1435.25
$ws.Range("J61").Value = 2236
$ws.Range("K61").Value = 1435.25
$ws.Range("L61").Value = 2236
$ws.Range("M61").Value = -1233.25
$ws.Range("N61").Value = -2640
$ws.Range("H68").Value = 2559.9048
$ws.Range("I68").Value = 2186.3572
$ws.Range("K68").Value = 2186.3572
$ws.Range("M68").Value = -1437.3572
$ws.Range("H71").Value = 2559.9048
$ws.Range("I71").Value = 2186.3572
$ws.Range("K71").Value = 10931.786
$ws.Range("M71").Value = -7187.786
$ws.Range("H113").Value = 1743.2307
$ws.Range("I113").Value = 1435.25
$ws.Range("J113").Value = 2236
$ws.Range("K113").Value = 1435.25
$ws.Range("L113").Value = 2236
$ws.Range("M113").Value = 734.75
$ws.Range("N113").Value = -6576

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2050.9
$ws.Range("I126").Value = 2282.0667
$ws.Range("K126").Value = 6846.2001
$ws.Range("M126").Value = -4376.2001
$ws.Range("H136").Value = 1333.12
$ws.Range("I136").Value = 920.91113
$ws.Range("J136").Value = 5043
$ws.Range("K136").Value = 2762.73339
$ws.Range("L136").Value = 15129
$ws.Range("M136").Value = -212.7333899999999
$ws.Range("N136").Value = -20229
